$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / date changed) for rows 2 through 10
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09)
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
